$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3073
$ws1.Range("F3").Value = 485
$ws1.Range("F4").Value = 59
$ws1.Range("F5").Value = 48
$ws1.Range("F6").Value = 3
$ws1.Range("F8").Value = 8
$ws1.Range("F9").Value = 1056
$ws1.Range("F10").Value = 14889
$ws1.Range("F11").Value = 182
$ws1.Range("F13").Value = 398
$ws1.Range("F14").Value = 5926
$ws1.Range("F17").Value = 52
$ws1.Range("F18").Value = 87
$ws1.Range("F19").Value = 1245
$ws1.Range("F20").Value = 20
$ws1.Range("F21").Value = 96
$ws1.Range("F22").Value = 197
$ws1.Range("F23").Value = 815
$ws1.Range("F24").Value = 2956
$ws1.Range("F26").Value = 10744
$ws1.Range("F27").Value = 1213
$ws1.Range("F28").Value = 81
$ws1.Range("F29").Value = 121
$ws1.Range("F30").Value = 3756
$ws1.Range("F31").Value = 251

# Sheet "全部类型" (sheet4): row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3073
$ws4.Range("F4").Value = 485
$ws4.Range("F5").Value = 59
$ws4.Range("F6").Value = 48
$ws4.Range("F7").Value = 3
$ws4.Range("F9").Value = 8
$ws4.Range("F10").Value = 1056
$ws4.Range("F11").Value = 14889
$ws4.Range("F12").Value = 182
$ws4.Range("F14").Value = 399
$ws4.Range("F15").Value = 5926
$ws4.Range("F18").Value = 52
$ws4.Range("F19").Value = 87
$ws4.Range("F20").Value = 1245
$ws4.Range("F21").Value = 20
$ws4.Range("F22").Value = 96
$ws4.Range("F23").Value = 197
$ws4.Range("F24").Value = 815
$ws4.Range("F25").Value = 2956
$ws4.Range("F28").Value = 10744
$ws4.Range("F29").Value = 1213
$ws4.Range("F30").Value = 81
$ws4.Range("F31").Value = 121
$ws4.Range("F32").Value = 3756
$ws4.Range("F33").Value = 251
